# Apply updated crypto price/volume figures (GitHub Actions scheduled refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.107.20'
$ws.Range("E2").Value = '  +4.09%  '
$ws.Range("D3").Value = '3.232.94'
$ws.Range("E3").Value = '  +2.15%  '
$ws.Range("D5").Value = '''576.50'
$ws.Range("E5").Value = '  +2.07%  '
$ws.Range("D6").Value = '''180.15'
$ws.Range("E6").Value = '  +5.49%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("E8").Value = '  -3.53%  '
$ws.Range("D9").Value = '3.232.30'
$ws.Range("E9").Value = '  +2.33%  '
$ws.Range("E10").Value = '  +4.01%  '
$ws.Range("E11").Value = '  +3.33%  '
$ws.Range("D12").Value = '''0.412'
$ws.Range("E12").Value = '  +4.45%  '
$ws.Range("D13").Value = '3.795.22'
$ws.Range("E13").Value = '  +2.40%  '
$ws.Range("E14").Value = '  +1.16%  '
$ws.Range("D15").Value = '''27.84'
$ws.Range("E15").Value = '  +2.02%  '
$ws.Range("D16").Value = '67.091.94'
$ws.Range("E16").Value = '  +4.30%  '
$ws.Range("E17").Value = '  +2.65%  '
$ws.Range("D18").Value = '3.246.75'
$ws.Range("E18").Value = '  +2.86%  '
$ws.Range("D19").Value = '''5.79'
$ws.Range("E19").Value = '  +0.89%  '
$ws.Range("D20").Value = '''13.36'
$ws.Range("E20").Value = '  +3.10%  '
$ws.Range("D21").Value = '''372.67'
$ws.Range("E21").Value = '  +4.95%  '
$ws.Range("D22").Value = '''7.55'
$ws.Range("E22").Value = '  +4.35%  '
$ws.Range("D23").Value = '''0.999'
$ws.Range("E23").Value = '  -0.36%  '
$ws.Range("D24").Value = '''70.73'
$ws.Range("E24").Value = '  +3.82%  '
$ws.Range("D25").Value = '''0.508'
$ws.Range("E25").Value = '  +1.47%  '
$ws.Range("E26").Value = '  +1.08%  '
$ws.Range("E27").Value = '  -0.44%  '
$ws.Range("E28").Value = '  +2.91%  '
$ws.Range("E29").Value = '  +0.23%  '
$ws.Range("E30").Value = '  +4.00%  '
$ws.Range("E31").Value = '  +4.88%  '
$ws.Range("D32").Value = '''22.53'
$ws.Range("E32").Value = '  +2.46%  '
$ws.Range("E34").Value = '  +4.44%  '
$ws.Range("D35").Value = '''6.84'
$ws.Range("E35").Value = '  +2.83%  '
$ws.Range("B36").Value = 'Monero'
$ws.Range("C36").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D36").Value = '''162.12'
$ws.Range("E36").Value = '  +5.41%  '
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").Value = '''1.49'
$ws.Range("E37").Value = '  +3.70%  '
$ws.Range("D38").Value = '''0.855'
$ws.Range("E38").Value = '  +3.77%  '
$ws.Range("E39").Value = '  +7.82%  '
$ws.Range("D40").Value = '''6.77'
$ws.Range("E40").Value = '  +13.41%  '
$ws.Range("D41").Value = '''26.64'
$ws.Range("E41").Value = '  +1.80%  '
$ws.Range("E42").Value = '  +3.59%  '
$ws.Range("D43").Value = '''360.29'
$ws.Range("E43").Value = '  +12.93%  '
$ws.Range("D44").Value = '''4.39'
$ws.Range("E44").Value = '  +5.19%  '
$ws.Range("D45").Value = '2.698.77'
$ws.Range("E45").Value = '  +1.93%  '
$ws.Range("D46").Value = '''25.51'
$ws.Range("E46").Value = '  +5.90%  '
$ws.Range("D47").Value = '''40.46'
$ws.Range("D48").Value = '''0.0671'
$ws.Range("E48").Value = '  +2.82%  '
$ws.Range("E49").Value = '  +1.80%  '
$ws.Range("B50").Value = 'Stellar'
$ws.Range("C50").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D50").Value = '''0.102'
$ws.Range("E50").Value = '  +1.29%  '
$ws.Range("B51").Value = 'ONDO'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D51").Value = '''0.993'
$ws.Range("E51").Value = '  +5.85%  '
